$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two email cells (C2/C3) -----------------------------------
# NOTE: C3 is written first on purpose. The workbook's shared-string table
# assigns indices in write order for newly-seen unique strings; writing C3
# ("abh...") before C2 ("nityaranjn...") reproduces the same shared-string
# ordering (and therefore the same <v> index references) as the target file.
$ws.Range("C3").Value = "abh0906536897454546@gmail.com"
$ws.Range("C2").Value = "nityaranjn78977344643@gmail.com"

# --- Widen column C (was 29.88671875 "best fit" -> 32.88671875) -----------
$ws.Columns("C").ColumnWidth = 32

# --- Move the active selection from C13 to C6 ------------------------------
$ws.Range("C6").Select()

# --- Resize the workbook window (best effort; persisted window chrome) ----
$win = $wb.Windows.Item(1)
$win.Width = 12864
$win.Height = 5448
